# corrected depression and reanalyzing it
#
# The patient_log sheet had several rows with a missing "Time_LHC" (column E)
# value even though the rest of the row had already been filled in, and the
# most recently logged patients (IDs 07027 through 07063) were missing their
# Date / Start_ECG / Time_LHC values altogether. This fills in the values
# from the re-analysis.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Fill in missing Time_LHC (column E) values for rows that already had data
# ---------------------------------------------------------------------------
$timeLhcUpdates = @(
    ,("E6",  0.58333333333333337)
    ,("E7",  0.61319444444444449)
    ,("E8",  0.46388888888888885)
    ,("E9",  0.39652777777777781)
    ,("E10", 0.64861111111111114)
    ,("E11", 0.41250000000000003)
    ,("E12", 0.44375000000000003)
    ,("E13", 0.53888888888888886)
    ,("E14", 0.57500000000000007)
    ,("E23", 0.73055555555555562)
    ,("E24", 0.49513888888888885)
    ,("E25", 0.46805555555555550)
    ,("E26", 0.58124999999999993)
    ,("E27", 0.53263888888888888)
    ,("E28", 0.47083333333333338)
    ,("E29", 0.51388888888888895)
    ,("E30", 0.48680555555555555)
)

foreach ($pair in $timeLhcUpdates) {
    $addr = $pair[0]
    $value = $pair[1]
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "h:mm"
    $cell.Value = $value
}

# ---------------------------------------------------------------------------
# Fill in missing Start_ECG (column C) values for rows 23, 24 and 27-30
# ---------------------------------------------------------------------------
$startEcgUpdates = @(
    ,("C23", 0.43611111111111112)
    ,("C24", 0.39513888888888887)
    ,("C27", 0.39652777777777781)
    ,("C28", 0.40069444444444446)
    ,("C29", 0.44444444444444442)
    ,("C30", 0.45833333333333331)
)

foreach ($pair in $startEcgUpdates) {
    $addr = $pair[0]
    $value = $pair[1]
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "h:mm"
    $cell.Value = $value
}

# ---------------------------------------------------------------------------
# Fill in missing Date (column B) values for rows 23-30
# ---------------------------------------------------------------------------
$dateUpdates = @(
    ,("B23", "10/14/2019")
    ,("B24", "10/17/2019")
    ,("B25", "10/21/2019")
    ,("B26", "10/21/2019")
    ,("B27", "10/24/2019")
    ,("B28", "10/30/2019")
    ,("B29", "10/30/2019")
    ,("B30", "11/11/2019")
)

foreach ($pair in $dateUpdates) {
    $addr = $pair[0]
    $value = $pair[1]
    $ws.Range($addr).Value = $value
}

# Leave the active selection on E15, matching the state the workbook was
# saved in after the edits were made.
$ws.Range("E15").Select() | Out-Null
